$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new "Season" column to the left of the table, and a new row
# below the existing data row for the 23/24 season entry.
$ws.Columns.Item(1).Insert()
$ws.Rows.Item(3).Insert()

# The table's underlying range moved with the insert but the ListObject
# definition needs to be told about its new extent explicitly.
$lo.Resize($ws.Range("B1:N3"))

# New "Season" column values
$ws.Range("A1").Value = "Season"
$ws.Range("A2").Value = "24/25"
$ws.Range("A3").Value = "23/24"

# Brighton's 23/24 season stats (new row)
$ws.Range("B3").Value = "Brighton"
$ws.Range("C3").Value = 11
$ws.Range("D3").Value = 48
$ws.Range("E3").Value = 38
$ws.Range("F3").Value = 12
$ws.Range("G3").Value = 14
$ws.Range("H3").Value = 12
$ws.Range("I3").Value = 55
$ws.Range("J3").Value = 62
$ws.Range("K3").Value = -7
$ws.Range("L3").Value = "João Pedro"
$ws.Range("M3").Value = 3
$ws.Range("N3").Value = 91

# Turn on a (worksheet-level) AutoFilter anchored at the new Season header
# cell, and register the hidden _FilterDatabase defined name Excel keeps
# alongside it.
$ws.Range("A1").AutoFilter() | Out-Null
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "='Brighton Stats'!`$A`$1:`$A`$1")
$filterName.Visible = $false

# Match the author's final selection
$ws.Range("G10").Select() | Out-Null
